$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "002/DR002/AV1"
$ws.Range("C2").Value = "BJ49785"
$ws.Range("D2").Value = "Anas tawfiqi"
$ws.Range("G2").Value = 15
$ws.Range("I2").Value = 10000
$ws.Range("K2").Value = 1500
$ws.Range("L2").Value = 7500
$ws.Range("M2").Value = 16000

# Row 3
$ws.Range("A3").Value = "002/DR002/AV1"
$ws.Range("C3").Value = "BB779645"
$ws.Range("D3").Value = "Karami abdelilah"
$ws.Range("I3").Value = 7500
$ws.Range("K3").Value = 750
$ws.Range("L3").Value = 7500
$ws.Range("M3").Value = 14250

# Row 4
$ws.Range("A4").Value = "108/ANSYSFYSN01/AV1"
$ws.Range("C4").Value = "'110384"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "real madrid fc"
$ws.Range("H4").Value = 10000
$ws.Range("I4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 10000

# Row 5
$ws.Range("A5").Value = "108/ANSYSFYSN01/AV1"
$ws.Range("C5").Value = "KS10293"
$ws.Range("D5").Value = "Karim benzima"
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 3000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 300
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 2700

# Row 6
$ws.Range("A6").Value = "108/ANSYSFYSN01/AV1"
$ws.Range("B6").Value = "Direction régionale"
$ws.Range("C6").Value = "BB102938"
$ws.Range("D6").Value = "Rodrigo silva jr"
$ws.Range("E6").Value = "non"
$ws.Range("H6").Value = 2000
$ws.Range("I6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 2000

# Row 7
$ws.Range("A7").Value = "108/ANSYSFYSN01/AV1"
$ws.Range("B7").Value = "Direction régionale"
$ws.Range("C7").Value = "BJ119649"
$ws.Range("D7").Value = "Test Test"
$ws.Range("G7").Value = 10
$ws.Range("H7").Value = 4000
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = 3600

# Row 8 (totals)
$ws.Range("H8").Value = 19000
$ws.Range("I8").Value = 17500
$ws.Range("J8").Value = 700
$ws.Range("K8").Value = 2250
$ws.Range("L8").Value = 15000
$ws.Range("M8").Value = 48550
